$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'318.64"
$ws.Range("E2").Value = "'4.28%"
$ws.Range("D3").Value = "'36.09"
$ws.Range("E3").Value = "'0.00%"
$ws.Range("D4").Value = "'5.133"
$ws.Range("E4").Value = "'0.71%"
$ws.Range("D5").Value = "'0.08209"
$ws.Range("E5").Value = "'4.42%"
$ws.Range("E6").Value = "'-1.19%"
$ws.Range("D7").Value = "'8.009"
$ws.Range("E7").Value = "'1.11%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9259"
$ws.Range("E8").Value = "'0.78%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1006"
$ws.Range("E9").Value = "'3.75%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1897"
$ws.Range("E10").Value = "'1.61%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.09220"
$ws.Range("E11").Value = "'6.37%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03598"
$ws.Range("E12").Value = "'3.14%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09920"
$ws.Range("E13").Value = "'-0.13%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001432"
$ws.Range("E14").Value = "'-0.37%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005677"
$ws.Range("E15").Value = "'-0.62%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.464"
$ws.Range("E16").Value = "'0.21%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'4.134"
$ws.Range("E17").Value = "'0.81%"
$ws.Range("E18").Value = "'16.82%"
$ws.Range("D19").Value = "'0.3374"
$ws.Range("E19").Value = "'-1.48%"
$ws.Range("E20").Value = "'2.25%"
$ws.Range("D21").Value = "'5.062"
$ws.Range("E21").Value = "'4.06%"
$ws.Range("D22").Value = "'0.2189"
$ws.Range("E22").Value = "'-0.57%"
$ws.Range("D23").Value = "'0.04590"
$ws.Range("E23").Value = "'0.90%"
$ws.Range("D24").Value = "'0.001243"
$ws.Range("E24").Value = "'0.98%"
$ws.Range("D25").Value = "'0.004736"
$ws.Range("E25").Value = "'-6.79%"
$ws.Range("D26").Value = "'0.0001302"
$ws.Range("E26").Value = "'-7.04%"
$ws.Range("D27").Value = "'0.0004501"
$ws.Range("E27").Value = "'-5.29%"
$ws.Range("D39").Value = "'0.02010"
$ws.Range("E39").Value = "'9.93%"
$ws.Range("D40").Value = "'0.04989"
$ws.Range("E40").Value = "'4.48%"
$ws.Range("D41").Value = "'0.007727"
$ws.Range("E41").Value = "'1.19%"
$ws.Range("D42").Value = "'0.1400"
$ws.Range("E42").Value = "'0.14%"
$ws.Range("D43").Value = "'0.007802"
$ws.Range("E43").Value = "'0.70%"
$ws.Range("D44").Value = "'0.002133"
$ws.Range("E44").Value = "'-4.38%"
$ws.Range("D45").Value = "'0.01200"
$ws.Range("E45").Value = "'5.97%"
$ws.Range("E46").Value = "'0.77%"
$ws.Range("E47").Value = "'-0.03%"
$ws.Range("E48").Value = "'18.29%"
$ws.Range("D49").Value = "'0.001901"
$ws.Range("E49").Value = "'-5.02%"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'-0.03%"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'-0.03%"
